$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column H ("type") for rows 3-8 (tags a2..a7) from FLOAT to INT
$ws.Range("H3").Value = "INT"
$ws.Range("H4").Value = "INT"
$ws.Range("H5").Value = "INT"
$ws.Range("H6").Value = "INT"
$ws.Range("H7").Value = "INT"
$ws.Range("H8").Value = "INT"

# Reset column F ("address") to 0 for all data rows (2-16)
for ($r = 2; $r -le 16; $r++) {
    $ws.Cells.Item($r, 6).Value = 0
}

# Unify body cell style: make rows 2-16, columns A-I use the same style as the
# current "A2" style (drop the distinct H/I FLOAT/DCBA look) by copying A2's style.
$bodyStyleSource = $ws.Range("A2")
for ($r = 2; $r -le 16; $r++) {
    $rowRange = $ws.Range("A" + $r + ":I" + $r)
    $rowRange.Style = $bodyStyleSource.Style
}

# Unify header style: make H1/I1 use the same style as A1
$ws.Range("H1:I1").Style = $ws.Range("A1").Style

# Update the selection to H7
$ws.Range("H7").Select()
